$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20; existing rows 20..68 shift down to 21..69.
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with its data.
$ws.Cells.Item(20, 1).Value()  = 11
$ws.Cells.Item(20, 2).Value()  = "Vega Monumental Concepción"
$ws.Cells.Item(20, 3).Value()  = "Bíobío"
$ws.Cells.Item(20, 4).Value()  = 44804
$ws.Cells.Item(20, 5).Value()  = 8
$ws.Cells.Item(20, 6).Value()  = 100112037
$ws.Cells.Item(20, 7).Value()  = "Cebollín"
$ws.Cells.Item(20, 8).Value()  = "Sin especificar"
$ws.Cells.Item(20, 9).Value()  = "Primera"
$ws.Cells.Item(20, 10).Value() = 150
$ws.Cells.Item(20, 11).Value() = 5000
$ws.Cells.Item(20, 12).Value() = 5500
$ws.Cells.Item(20, 13).Value() = 5333
$ws.Cells.Item(20, 14).Value() = "$/paquete 36 unidades"
$ws.Cells.Item(20, 15).Value() = "Región Metropolitana"
$ws.Cells.Item(20, 16).Value() = 148
$ws.Cells.Item(20, 17).Value() = 36
$ws.Cells.Item(20, 18).Value() = "Hortaliza"
